$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 159. Excel shifts the existing
# rows 159-199 down to 160-200 (and the sheet's used-range / <dimension>
# grows from R199 to R200 automatically).
$ws.Rows.Item(159).Insert()

# Populate the newly inserted row 159 with the new Berenjena record.
$ws.Range("A159").Value = 10
$ws.Range("B159").Value = "Vega Modelo de Temuco"
$ws.Range("C159").Value = "La Araucanía"
$ws.Range("D159").Value = 44508
$ws.Range("E159").Value = 9
$ws.Range("F159").Value = 100112001
$ws.Range("G159").Value = "Berenjena"
$ws.Range("H159").Value = "Sin especificar"
$ws.Range("I159").Value = "Primera"
$ws.Range("J159").Value = 170
$ws.Range("K159").Value = 10000
$ws.Range("L159").Value = 11000
$ws.Range("M159").Value = 10529
$ws.Range("N159").Value = "$/caja 60 unidades"
$ws.Range("O159").Value = "Región de Arica y Parinacota"
$ws.Range("P159").Value = 175
$ws.Range("Q159").Value = 60
$ws.Range("R159").Value = "Hortaliza"

# Make sure the date cell keeps the existing date number format (style
# carries over from the Insert(), but set it explicitly to be safe).
$ws.Range("D159").NumberFormat = $ws.Range("D160").NumberFormat
